# Apply updated market price / profit figures to each profession sheet
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 712.4211
$ws.Range("I53").Value = 984.6667
$ws.Range("J53").Value = 586.7692
$ws.Range("K53").Value = 984.6667
$ws.Range("L53").Value = 586.7692
$ws.Range("M53").Value = -347.6667
$ws.Range("N53").Value = -1860.7692
$ws.Range("H74").Value = 4480.9375
$ws.Range("I74").Value = 4345.769
$ws.Range("J74").Value = 5066.6665
$ws.Range("K74").Value = 4345.769
$ws.Range("L74").Value = 5066.6665
$ws.Range("M74").Value = -3409.769
$ws.Range("N74").Value = -6938.6665
$ws.Range("H77").Value = 4480.9375
$ws.Range("I77").Value = 4345.769
$ws.Range("J77").Value = 5066.6665
$ws.Range("K77").Value = 21728.845
$ws.Range("L77").Value = 25333.3325
$ws.Range("M77").Value = -17048.845
$ws.Range("N77").Value = -34693.3325
$ws.Range("H80").Value = 5756.4546
$ws.Range("I80").Value = 290
$ws.Range("K80").Value = 870
$ws.Range("M80").Value = 128
$ws.Range("H83").Value = 5756.4546
$ws.Range("I83").Value = 290
$ws.Range("K83").Value = 2610
$ws.Range("M83").Value = 2382
$ws.Range("H107").Value = 2837.5
$ws.Range("I107").Value = 2837.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2837.5
$ws.Range("L107").Value = 0
$ws.Range("N107").Value = -917.5
$ws.Range("M107").ClearContents()
$ws.Range("H112").Value = 2020.421
$ws.Range("J112").Value = 2020.421
$ws.Range("L112").Value = 6061.263
$ws.Range("N112").Value = -8277.262999999999
$ws.Range("H125").Value = 922.2222
$ws.Range("I125").Value = 825
$ws.Range("J125").Value = 1000
$ws.Range("K125").Value = 7425
$ws.Range("L125").Value = 9000
$ws.Range("M125").Value = -4965
$ws.Range("N125").Value = -13920
$ws.Range("H138").Value = 2926.2
$ws.Range("J138").Value = 3699
$ws.Range("L138").Value = 11097
$ws.Range("N138").Value = -21377

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4503.1113
$ws.Range("I32").Value = 1893.569
$ws.Range("K32").Value = 1893.569
$ws.Range("M32").Value = -1606.569
$ws.Range("H39").Value = 10499
$ws.Range("I39").Value = 10998.75
$ws.Range("J39").Value = 8500
$ws.Range("K39").Value = 10998.75
$ws.Range("L39").Value = 8500
$ws.Range("M39").Value = -10478.75
$ws.Range("N39").Value = -9540
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H45").Value = 7881092.5
$ws.Range("I45").Value = 1749.25
$ws.Range("K45").Value = 1749.25
$ws.Range("M45").Value = -1372.25
$ws.Range("H47").Value = 10000
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H88").Value = 234000
$ws.Range("J88").Value = 234000
$ws.Range("L88").Value = 234000
$ws.Range("N88").Value = -234812
$ws.Range("H91").Value = 234000
$ws.Range("J91").Value = 234000
$ws.Range("L91").Value = 234000
$ws.Range("N91").Value = -236808
$ws.Range("H122").Value = 3851.5806
$ws.Range("I122").Value = 3541.0386
$ws.Range("K122").Value = 10623.1158
$ws.Range("M122").Value = -8173.1158
$ws.Range("H132").Value = 2981.1345
$ws.Range("I132").Value = 2600.7727
$ws.Range("J132").Value = 5073.125
$ws.Range("K132").Value = 7802.3181
$ws.Range("L132").Value = 15219.375
$ws.Range("M132").Value = -5272.3181
$ws.Range("N132").Value = -20279.375
$ws.Range("H137").Value = 89191.42999999999
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6103.625
$ws.Range("I86").Value = 3141
$ws.Range("J86").Value = 7450.273
$ws.Range("K86").Value = 3141
$ws.Range("L86").Value = 7450.273
$ws.Range("M86").Value = -2018
$ws.Range("N86").Value = -9696.273000000001
$ws.Range("H89").Value = 6103.625
$ws.Range("I89").Value = 3141
$ws.Range("J89").Value = 7450.273
$ws.Range("K89").Value = 15705
$ws.Range("L89").Value = 37251.365
$ws.Range("M89").Value = -10089
$ws.Range("N89").Value = -48483.365
$ws.Range("H105").Value = 75021.21000000001
$ws.Range("I105").Value = 80542.08
$ws.Range("K105").Value = 80542.08
$ws.Range("M105").Value = -78795.08

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 74.34999999999999
$ws.Range("I12").Value = 55.125
$ws.Range("J12").Value = 87.166664
$ws.Range("K12").Value = 165.375
$ws.Range("L12").Value = 261.499992
$ws.Range("M12").Value = 7.625
$ws.Range("N12").Value = -607.499992
$ws.Range("H86").Value = 553.2143
$ws.Range("I86").Value = 364.14285
$ws.Range("J86").Value = 742.2857
$ws.Range("K86").Value = 1092.42855
$ws.Range("L86").Value = 2226.8571
$ws.Range("M86").Value = 93.57144999999991
$ws.Range("N86").Value = -4598.8571
$ws.Range("H89").Value = 553.2143
$ws.Range("I89").Value = 364.14285
$ws.Range("J89").Value = 742.2857
$ws.Range("K89").Value = 3277.28565
$ws.Range("L89").Value = 6680.571300000001
$ws.Range("M89").Value = 2650.71435
$ws.Range("N89").Value = -18536.5713
$ws.Range("H137").Value = 3427.1765
$ws.Range("I137").Value = 2302.4285
$ws.Range("K137").Value = 6907.2855
$ws.Range("M137").Value = -1807.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6870
$ws.Range("I80").Value = 3750
$ws.Range("K80").Value = 3750
$ws.Range("M80").Value = -2752
$ws.Range("H83").Value = 6870
$ws.Range("I83").Value = 3750
$ws.Range("K83").Value = 18750
$ws.Range("M83").Value = -13758
$ws.Range("H97").Value = 2280.0715
$ws.Range("I97").Value = 2274.6365
$ws.Range("K97").Value = 2274.6365
$ws.Range("M97").Value = -1778.6365
$ws.Range("H119").Value = 77134.17999999999
$ws.Range("J119").Value = 77134.17999999999
$ws.Range("L119").Value = 77134.17999999999
$ws.Range("N119").Value = -86810.17999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6645.346
$ws.Range("J7").Value = 2931.25
$ws.Range("L7").Value = 2931.25
$ws.Range("N7").Value = -3155.25
$ws.Range("H22").Value = 159257.14
$ws.Range("I22").Value = 3249.5
$ws.Range("J22").Value = 221660.2
$ws.Range("K22").Value = 3249.5
$ws.Range("L22").Value = 221660.2
$ws.Range("M22").Value = -2954.5
$ws.Range("N22").Value = -222250.2
$ws.Range("H27").Value = 159257.14
$ws.Range("I27").Value = 3249.5
$ws.Range("J27").Value = 221660.2
$ws.Range("K27").Value = 3249.5
$ws.Range("L27").Value = 221660.2
$ws.Range("M27").Value = -3142.5
$ws.Range("N27").Value = -221874.2
$ws.Range("H35").Value = 543.2
$ws.Range("I35").Value = 420.25
$ws.Range("J35").Value = 1035
$ws.Range("K35").Value = 420.25
$ws.Range("L35").Value = 1035
$ws.Range("M35").Value = -84.25
$ws.Range("N35").Value = -1707
$ws.Range("H40").Value = 1953403.5
$ws.Range("I40").Value = 41698.04
$ws.Range("K40").Value = 41698.04
$ws.Range("M40").Value = -41562.04
$ws.Range("H46").Value = 21695.4
$ws.Range("I46").Value = 26744.25
$ws.Range("J46").Value = 1500
$ws.Range("K46").Value = 26744.25
$ws.Range("L46").Value = 1500
$ws.Range("M46").Value = -26556.25
$ws.Range("N46").Value = -1876
$ws.Range("H55").Value = 8222.933999999999
$ws.Range("I55").Value = 1956.25
$ws.Range("J55").Value = 15384.857
$ws.Range("K55").Value = 1956.25
$ws.Range("L55").Value = 15384.857
$ws.Range("M55").Value = -1783.25
$ws.Range("N55").Value = -15730.857
$ws.Range("H100").Value = 11797.706
$ws.Range("I100").Value = 15824.875
$ws.Range("J100").Value = 8218
$ws.Range("K100").Value = 15824.875
$ws.Range("L100").Value = 8218
$ws.Range("M100").Value = -15283.875
$ws.Range("N100").Value = -9300
$ws.Range("H122").Value = 77081040
$ws.Range("I122").Value = 100204136
$ws.Range("J122").Value = 4035.6667
$ws.Range("K122").Value = 300612408
$ws.Range("L122").Value = 12107.0001
$ws.Range("M122").Value = -300609958
$ws.Range("N122").Value = -17007.0001
$ws.Range("H126").Value = 6645.346
$ws.Range("J126").Value = 2931.25
$ws.Range("L126").Value = 8793.75
$ws.Range("N126").Value = -13733.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 48108.555
$ws.Range("I81").Value = 3853.5715
$ws.Range("J81").Value = 203001
$ws.Range("K81").Value = 7707.143
$ws.Range("L81").Value = 406002
$ws.Range("M81").Value = -6646.143
$ws.Range("N81").Value = -408124
$ws.Range("H84").Value = 48108.555
$ws.Range("I84").Value = 3853.5715
$ws.Range("J84").Value = 203001
$ws.Range("K84").Value = 38535.715
$ws.Range("L84").Value = 2030010
$ws.Range("M84").Value = -33231.715
$ws.Range("N84").Value = -2040618
$ws.Range("H100").Value = 3760270
$ws.Range("I100").Value = 4464971
$ws.Range("K100").Value = 8929942
$ws.Range("M100").Value = -8929401
$ws.Range("H122").Value = 2540.8572
$ws.Range("I122").Value = 2501.8333
$ws.Range("K122").Value = 7505.499899999999
$ws.Range("M122").Value = -5055.499899999999
$ws.Range("H132").Value = 2814.3242
$ws.Range("I132").Value = 2742.4075
$ws.Range("K132").Value = 8227.2225
$ws.Range("M132").Value = -5697.2225

Write-Output "Updated market price data across ALC, ARM, BSM, CUL, GSM, LTW, WVR sheets"